$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("org_list")

# The location "Le Sud-Ouest" was stored with a non-breaking space between
# "Le" and "Sud-Ouest" which broke downstream matching/localisation logic.
# Re-write the affected cells with a normal space so a new shared string
# is created and used instead of the old (non-breaking-space) one.
$fixed = "Le Sud-Ouest"

$ws.Range("C9").Value = $fixed
$ws.Range("C19").Value = $fixed
$ws.Range("C22").Value = $fixed

# Match the author's final selection location in the sheet.
$ws.Range("C25").Select()
